$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Bug #6 (row 6): "Doesn't work well when invoking from PDF highlighted word"
# mark as resolved -> strikethrough Priority + Description
$ws.Range("A6").Font.Strikethrough = $true
$ws.Range("B6").Font.Strikethrough = $true

# --- Bug #11 (row 11): "Opening new window for highlighted word ... doesn't work"
# mark as resolved -> strikethrough Priority + Description, add Solution text
$ws.Range("A11").Font.Strikethrough = $true
$ws.Range("B11").Font.Strikethrough = $true
$ws.Range("D11").Value = "Needed Application.DoEvents before pasting"
$ws.Range("D11").Font.Strikethrough = $false
$ws.Rows.Item(11).RowHeight = 30

# --- Bug #25 (row 25): "Let the WebBrowser control lazy-load ..."
# mark as resolved -> strikethrough Priority + Description, add Solution text
$ws.Range("A25").Font.Strikethrough = $true
$ws.Range("B25").Font.Strikethrough = $true
$ws.Range("D25").Value = "Changed control to use Url property instead of manually setting a static HTML content. Also  " + [char]10 + "control was previously using older IE versions to browse web."
$ws.Range("D25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 45

# Column D (Solution) is now much wider because of the new, longer solution text
$ws.Columns.Item(4).ColumnWidth = 78.59244791666667

# Last-used selection cell
[void]$ws.Range("B24").Select()
